$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the next batch of vocabulary rows (28-32) as "Processed" in column C,
# matching the other already-processed rows above them.
foreach ($r in 28..32) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
